$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("G6").Value = 2.25
$ws.Range("H6").Value = 3.2
$ws.Range("I6").Value = 3
$ws.Range("J6").Value = 2.82
$ws.Range("K6").Value = 2.05
$ws.Range("L6").Value = 3.6
$ws.Range("M6").Value = 1.02
$ws.Range("N6").Value = 9.65
$ws.Range("O6").Value = 1.28
$ws.Range("P6").Value = 3.05
$ws.Range("Q6").Value = 1.88
$ws.Range("R6").Value = 1.83
$ws.Range("S6").Value = 1.4
$ws.Range("T6").Value = 2.52
$ws.Range("U6").Value = 1.65
$ws.Range("V6").Value = 1.98
$ws.Range("W6").Value = 8.25
$ws.Range("X6").Value = 11.75
$ws.Range("Y6").Value = 8.75
$ws.Range("Z6").Value = 23
$ws.Range("AA6").Value = 17.5
$ws.Range("AB6").Value = 26
$ws.Range("AC6").Value = 9.75
$ws.Range("AD6").Value = 6.2
$ws.Range("AE6").Value = 13
$ws.Range("AF6").Value = 55
$ws.Range("AG6").Value = 9.25
$ws.Range("AH6").Value = 15.5
$ws.Range("AI6").Value = 10.75
$ws.Range("AJ6").Value = 37
$ws.Range("AK6").Value = 26
$ws.Range("AL6").Value = 32
$ws.Range("AM6").Value = 400
$ws.Range("AN6").Value = 4.15
$ws.Range("AO6").Value = 11.75
$ws.Range("AQ6").Value = 45
$ws.Range("AT6").Value = 2.5
$ws.Range("AU6").Value = 6.8
$ws.Range("AV6").Value = 60
$ws.Range("AW6").Value = 4.9
$ws.Range("AX6").Value = 17
$ws.Range("AY6").Value = 24
$ws.Range("AZ6").Value = 80
$ws.Range("BA6").Value = 120
$ws.Range("BB6").Value = 300

# Row 8
$ws.Range("O8").Value = 1.26
$ws.Range("Q8").Value = 1.82
$ws.Range("Z8").Value = 60
$ws.Range("AA8").Value = 35
$ws.Range("AD8").Value = 6.7
$ws.Range("AU8").Value = 7
$ws.Range("AX8").Value = 9.25
$ws.Range("AY8").Value = 17.5
$ws.Range("BB8").Value = 200

# Row 27
$ws.Range("J27").Value = 2.57
$ws.Range("L27").Value = 4.65
$ws.Range("W27").Value = 5.6
$ws.Range("AA27").Value = 18.5
$ws.Range("AB27").Value = 37
$ws.Range("AC27").Value = 6.4
$ws.Range("AD27").Value = 6
$ws.Range("AG27").Value = 8.5
$ws.Range("AL27").Value = 70
$ws.Range("AN27").Value = 3.6
$ws.Range("AO27").Value = 10
$ws.Range("AP27").Value = 21
$ws.Range("AR27").Value = 80
$ws.Range("AT27").Value = 2.22
$ws.Range("AU27").Value = 7.8
$ws.Range("AW27").Value = 5.7
$ws.Range("AY27").Value = 37
$ws.Range("AZ27").Value = 175
$ws.Range("BA27").Value = 250
